# StructureDefinition-canonical-measure.xlsx update
# Rebrand from "Alvearie"/ibm.com to "LinuxForHealth"/linuxforhealth.org,
# bump version 7.0.0 -> 8.0.0, update the publish date, and clear the
# stray ele-1/ext-1 constraint text that had been duplicated onto the
# root "Extension" row of the Elements sheet.

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/canonical-measure"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# The "Fixed Value" for Extension.url mirrors the canonical URL above.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/canonical-measure"

# The root "Extension" row no longer carries the ele-1/ext-1 constraint
# text in its "Constraint(s)" column (it still appears, correctly, on
# the Extension.extension row below).
$elements.Range("AI2").Value = ""
